$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    # 1 = wdReplaceOne (replace only the single match within this
    # paragraph-scoped range; the doc has no duplicate paragraph text so
    # this is both safe and sufficient).
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
    if (-not $found) {
        Write-Host "WARNING: text not found in paragraph $paraIndex -> $oldText"
    }
}

# "Objetivos" body paragraph becomes the old "Programa resumido" text
Replace-InParagraph 6 "Apresentar os princípios da automação da produção, características, aplicações e capacidades" "Controle e automação; Robótica; Domótica; Sistemas Supervisórios, Pneumática, Hidráulica, CLP"

# "Docente(s) Responsável(eis)" list item becomes the old "Objetivos" text
Replace-InParagraph 8 "5840917 - Fabricio Maciel Gomes" "Apresentar os princípios da automação da produção, características, aplicações e capacidades"

# "Programa resumido" body paragraph becomes the old "Programa" text
Replace-InParagraph 10 "Controle e automação; Robótica; Domótica; Sistemas Supervisórios, Pneumática, Hidráulica, CLP" "Introdução aos princípios de controle e automação; Fundamentos da Robótica; Fundamentos da Domótica;  Introdução a Sistemas Supervisórios, Princípios da Automação Pneumática, Hidráulica, Introdução aos Controladores Lógicos Programáveis."

# "Programa" body paragraph becomes the old "Método:" value text
Replace-InParagraph 12 "Introdução aos princípios de controle e automação; Fundamentos da Robótica; Fundamentos da Domótica;  Introdução a Sistemas Supervisórios, Princípios da Automação Pneumática, Hidráulica, Introdução aos Controladores Lógicos Programáveis." "Aulas expositivas e práticas."

# Paragraph 14 ("Avaliação" list item) holds three bold labels (Método /
# Critério / Norma de recuperação) each followed by a value run and a
# <w:br/> line break, all inside a single Word paragraph. Several of the
# new values coincide with other old values in this same paragraph, so
# plain text Find/Replace cannot reliably tell the runs apart once the
# first substitution lands. Instead, locate the three label runs (which
# stay unique) to carve out exact, disjoint Range objects for each value,
# then assign .Text on those ranges - last to first - so earlier offsets
# remain valid while later ones are still being computed.
function Get-LabelRange($searchRange, $labelText) {
    $r = $searchRange.Duplicate()
    $found = $r.Find.Execute($labelText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "label not found: $labelText"
    }
    return $r
}

$p14 = $d.Paragraphs(14)
$para14Start = $p14.Range.Start
$para14End = $p14.Range.End

$metodoLabel = Get-LabelRange ($d.Range($para14Start, $para14End)) "Método: "
$criterioLabel = Get-LabelRange ($d.Range($metodoLabel.End, $para14End)) "Critério: "
$normaLabel = Get-LabelRange ($d.Range($criterioLabel.End, $para14End)) "Norma de recuperação: "

$metodoValRange = $d.Range($metodoLabel.End, $criterioLabel.Start - 1)
$criterioValRange = $d.Range($criterioLabel.End, $normaLabel.Start - 1)
$normaValRange = $d.Range($normaLabel.End, $para14End - 1)

# Sanity-check we sliced the expected original text before overwriting it.
if ($metodoValRange.Text -ne "Aulas expositivas e práticas.") { throw "Metodo value range mismatch: $($metodoValRange.Text)" }
if ($criterioValRange.Text -ne "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2") { throw "Criterio value range mismatch: $($criterioValRange.Text)" }
if ($normaValRange.Text -ne "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação.") { throw "Norma value range mismatch: $($normaValRange.Text)" }

# Assign from the end of the paragraph backwards so earlier Range objects
# (computed above, before any mutation) keep pointing at the right offsets.
$normaValRange.Text = "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) NISE, N. S., “Engenharia de Sistemas de Controle”, 3ª ed., LTC, 2002. OGATA, K., “Engenharia de Controle Moderno”, 4ª ed., Prentice-Hall do Brasil, 2003. Tutoriais disponibilizados pelo professor BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U. B.. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p. CAPELLI, A. Automação Industrial: controle de movimento e processos contínuos. São Paulo: Érica, 2006. SILVEIRA, P. R. da; SANTOS, W. E. Automação e controle discreto. 3. ed. São Paulo: Érica, 1998. MORAES, C. C.; CATRUCCI, P. Engenharia de automação industrial. 2. ed. Rio de Janeiro: LTC, 2007. GIORGINI, M. Automação aplicada: descrição e implementação de sistemas sequencias com PLC's. 5. ed. São Paulo: Érica, 2003."
$criterioValRange.Text = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$metodoValRange.Text = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"

# "Bibliografia" body paragraph becomes the old "Docente(s)" professor entry
Replace-InParagraph 16 "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) NISE, N. S., “Engenharia de Sistemas de Controle”, 3ª ed., LTC, 2002. OGATA, K., “Engenharia de Controle Moderno”, 4ª ed., Prentice-Hall do Brasil, 2003. Tutoriais disponibilizados pelo professor BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U. B.. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p. CAPELLI, A. Automação Industrial: controle de movimento e processos contínuos. São Paulo: Érica, 2006. SILVEIRA, P. R. da; SANTOS, W. E. Automação e controle discreto. 3. ed. São Paulo: Érica, 1998. MORAES, C. C.; CATRUCCI, P. Engenharia de automação industrial. 2. ed. Rio de Janeiro: LTC, 2007. GIORGINI, M. Automação aplicada: descrição e implementação de sistemas sequencias com PLC's. 5. ed. São Paulo: Érica, 2003." "5840917 - Fabricio Maciel Gomes"

Write-Host "Done."
